# Build the list of cell updates (row -> {Column = NewValue}) taken from the
# latest cryptos.com scrape. Column D ("Price") values that are purely numeric
# need to stay text, so we force the cell's number format to Text ("@") before
# writing them -- otherwise Excel would silently convert e.g. "2.27" into the
# number 2.27 and drop the original text representation (thousand-dot prices
# like "64.416.68" are never parsed as numbers, so they don't need this).
$updates = @(
    @{ Cell = "D2"; Value = "64.416.68"; ForceText = $false }
    @{ Cell = "D3"; Value = "3.186.73"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +2.87%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  +0.11%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "594.19"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +1.74%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "148.78"; ForceText = $true }
    @{ Cell = "E6"; Value = "  +2.76%  "; ForceText = $false }
    @{ Cell = "E7"; Value = "  +0.08%  "; ForceText = $false }
    @{ Cell = "D8"; Value = "3.175.64"; ForceText = $false }
    @{ Cell = "E8"; Value = "  +2.71%  "; ForceText = $false }
    @{ Cell = "E9"; Value = "  +1.27%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "0.163"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +2.08%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "6.07"; ForceText = $true }
    @{ Cell = "E11"; Value = "  +7.97%  "; ForceText = $false }
    @{ Cell = "E12"; Value = "  +1.38%  "; ForceText = $false }
    @{ Cell = "E13"; Value = "  +1.26%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "37.90"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +2.39%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "3.713.89"; ForceText = $false }
    @{ Cell = "E15"; Value = "  +2.86%  "; ForceText = $false }
    @{ Cell = "E17"; Value = "  +4.34%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "3.180.83"; ForceText = $false }
    @{ Cell = "E18"; Value = "  +2.77%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "64.206.02"; ForceText = $false }
    @{ Cell = "E19"; Value = "  +1.34%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "472.53"; ForceText = $true }
    @{ Cell = "E20"; Value = "  +2.77%  "; ForceText = $false }
    @{ Cell = "E21"; Value = "  +2.78%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "0.741"; ForceText = $true }
    @{ Cell = "E22"; Value = "  +2.60%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "7.71"; ForceText = $true }
    @{ Cell = "E23"; Value = "  +3.91%  "; ForceText = $false }
    @{ Cell = "E24"; Value = "  +10.01%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "13.33"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +3.21%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "81.91"; ForceText = $true }
    @{ Cell = "E26"; Value = "  +1.09%  "; ForceText = $false }
    @{ Cell = "E27"; Value = "  +0.14%  "; ForceText = $false }
    @{ Cell = "E28"; Value = "  +8.72%  "; ForceText = $false }
    @{ Cell = "E29"; Value = "  +2.69%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "2.27"; ForceText = $true }
    @{ Cell = "E30"; Value = "  +2.95%  "; ForceText = $false }
    @{ Cell = "E31"; Value = "  +0.12%  "; ForceText = $false }
    @{ Cell = "E32"; Value = "  +4.48%  "; ForceText = $false }
    @{ Cell = "E33"; Value = "  +10.53%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "28.44"; ForceText = $true }
    @{ Cell = "E34"; Value = "  +6.82%  "; ForceText = $false }
    @{ Cell = "E35"; Value = "  +2.04%  "; ForceText = $false }
    @{ Cell = "E36"; Value = "  +4.18%  "; ForceText = $false }
    @{ Cell = "E37"; Value = "  +4.63%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "2.32"; ForceText = $true }
    @{ Cell = "E38"; Value = "  +0.92%  "; ForceText = $false }
    @{ Cell = "E39"; Value = "  +0.51%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "469.23"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +7.98%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "9.43"; ForceText = $true }
    @{ Cell = "E41"; Value = "  +8.72%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "51.55"; ForceText = $true }
    @{ Cell = "E42"; Value = "  +2.58%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "0.297"; ForceText = $true }
    @{ Cell = "E43"; Value = "  +8.16%  "; ForceText = $false }
    @{ Cell = "E44"; Value = "  +2.59%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "2.939.10"; ForceText = $false }
    @{ Cell = "E45"; Value = "  +1.98%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "39.49"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +8.67%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "133.04"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +5.93%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "2.28"; ForceText = $true }
    @{ Cell = "E50"; Value = "  +6.29%  "; ForceText = $false }
    @{ Cell = "E51"; Value = "  +1.54%  "; ForceText = $false }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    if ($update.ForceText) {
        $range.NumberFormat = "@"
    }
    $range.Value = $update.Value
}
